# Apply the "Tableau_RESF_2024_1" edit:
#  - remove rows 15-21 (old Prelevements-obligatoires-corrige / Dette publique block is
#    collapsed into the remaining 14-row table, and the stray duplicate rows disappear)
#  - rewrite row 1..14 content & turn the bare numeric figures into "xx,x%" strings
#  - shrink the used range to A1:D14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the trailing rows (15-21) first, so row numbers below stay stable ---
$ws.Range("A15:D21").EntireRow.Delete()

# --- 2. Row 1 : header becomes blank / "2023" / "2024" ---
# (leading apostrophe forces these year-looking values to stay TEXT, not numbers,
#  matching the source t="inlineStr" cells; a lone "'" yields an empty TEXT cell
#  instead of clearing the cell outright)
$ws.Cells.Item(1,1).Value = "'"
$ws.Cells.Item(1,2).Value = "'2023"
$ws.Cells.Item(1,3).Value = "'2024"

# --- 3. Row 2 : markdown-style separator ---
$ws.Cells.Item(2,1).Value = ":--"
$ws.Cells.Item(2,2).Value = ":--:"
$ws.Cells.Item(2,3).Value = ":--:"

# --- 4. Row 3 : PIB ---
$ws.Cells.Item(3,1).Value = "PIB (PLF 2024)"
$ws.Cells.Item(3,2).Value = "1,0%"
$ws.Cells.Item(3,3).Value = "1,4%"

# --- 5. Row 4 : Indice des prix ... ---
$ws.Cells.Item(4,1).Value = "Indice des prix à la consommation harmonisé (PLF 2024)"
$ws.Cells.Item(4,2).Value = "5,7%"
$ws.Cells.Item(4,3).Value = "2,7%"

# --- 6. Row 5 : Solde public ---
$ws.Cells.Item(5,1).Value = "Solde public (en % du PIB) (PLF 2024)"
$ws.Cells.Item(5,2).Value = "-4,9%"
$ws.Cells.Item(5,3).Value = "-4,4%"

# --- 7. Row 6 : now a blank label + year header row ---
$ws.Cells.Item(6,1).Value = "'"
$ws.Cells.Item(6,2).Value = "'2022"
$ws.Cells.Item(6,3).Value = "'2023"
$ws.Cells.Item(6,4).Value = "'2024"

# --- 8. Row 7 : markdown-style separator (now 4 columns) ---
$ws.Cells.Item(7,1).Value = ":--"
$ws.Cells.Item(7,2).Value = ":--:"
$ws.Cells.Item(7,3).Value = ":--:"
$ws.Cells.Item(7,4).Value = ":--:"

# --- 9. Row 8 : Depense publique hors credits d'impot ---
$ws.Cells.Item(8,1).Value = "Dépense publique hors crédits d'impôt (en % du PIB)"
$ws.Cells.Item(8,2).Value = "57,7%"
$ws.Cells.Item(8,3).Value = "55,9%"
$ws.Cells.Item(8,4).Value = "55,3%"

# --- 10. Row 9 : Progression en volume des depenses publiques ---
$ws.Cells.Item(9,1).Value = "Progression en volume des dépenses publiques"
$ws.Cells.Item(9,2).Value = "-1,1%"
$ws.Cells.Item(9,3).Value = "-1,3%"
$ws.Cells.Item(9,4).Value = "0,5%"

# --- 11. Row 10 : Prelevements obligatoires ---
$ws.Cells.Item(10,1).Value = "Prélèvements obligatoires (en % du PIB)"
$ws.Cells.Item(10,2).Value = "45,4%"
$ws.Cells.Item(10,3).Value = "44,0%"
$ws.Cells.Item(10,4).Value = "44,1%"

# --- 12. Row 11 : blank label + year header row ---
$ws.Cells.Item(11,1).Value = "'"
$ws.Cells.Item(11,2).Value = "'2022"
$ws.Cells.Item(11,3).Value = "'2023"
$ws.Cells.Item(11,4).Value = "'2024"

# --- 13. Row 12 : markdown-style separator (4 columns) ---
$ws.Cells.Item(12,1).Value = ":--"
$ws.Cells.Item(12,2).Value = ":--:"
$ws.Cells.Item(12,3).Value = ":--:"
$ws.Cells.Item(12,4).Value = ":--:"

# --- 14. Row 13 : Dette publique totale ---
$ws.Cells.Item(13,1).Value = "Dette publique totale (en % du PIB)"
$ws.Cells.Item(13,2).Value = "111,8%"
$ws.Cells.Item(13,3).Value = "109,7%"
$ws.Cells.Item(13,4).Value = "109,7%"

# --- 15. Row 14 : Dette publique hors soutien a la zone euro ---
$ws.Cells.Item(14,1).Value = "Dette publique hors soutien à la zone euro (en % du PIB)"
$ws.Cells.Item(14,2).Value = "109,3%"
$ws.Cells.Item(14,3).Value = "107,4%"
$ws.Cells.Item(14,4).Value = "107,6%"
